# This reverts an earlier "Update results (best so far)" commit: it
# removes the rows/cells that commit had added to Sheet1 and restores
# the sheet to its prior (smaller) extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# That commit had added two cells, H22 and I22, to row 22 - clear them
# (Clear, not ClearContents, so the cells disappear entirely along with
# their style, matching the earlier/target state).
$ws.Range("H22:I22").Clear()

# It also added a block of cells C23:N23 to row 23 - remove those too,
# leaving only the pre-existing A23/B23/O23 cells in that row.
$ws.Range("C23:N23").Clear()

# And it added an entirely new row 24 (A24/B24) - delete the whole row.
$ws.Rows("24:24").Delete()

# Restore the view/selection as closely as the host lets us: the sheet
# is no longer zoomed to 75%, and the last active cell was N25.
$excel.ActiveWindow.Zoom = 100
$ws.Range("N25").Select()

$wb.Save()
